$d = $word.ActiveDocument

# 1) Merge the split "Working (" / "Rao-Scott+F" / ") LRT for " runs (which
#    were wrapped with spell-check proofErr tags around "Rao-Scott+F") into
#    a single literal run of text "Working (Rao-Scott+F) LRT for ". This
#    phrase occurs twice in the document (Classical section and Causal
#    section), so replace all occurrences.
$d.Content.Find.Execute("Working (Rao-Scott+F) LRT for ", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Working (Rao-Scott+F) LRT for ", 2)

# 2) Update the statistic values in the "Working 2logLR" line (Causal section)
$d.Content.Find.Execute("Working 2logLR =  78.4 p= 0.0000000000002 ", $true, `
    $false, $false, $false, $false, $true, 1, $false, `
    "Working 2logLR =  87.3 p= 0.0000000000000008 ", 2)

# 3) Update the scale factors line (Causal section)
$d.Content.Find.Execute("(scale factors:  1.4 0.96 0.85 0.83 );  denominator ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "(scale factors:  1.2 1.1 0.86 0.81 );  denominator ", 2)
